$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Conferences")
$ws.Columns.Item(3).ColumnWidth = 12.140625
$ws.Columns.Item(4).ColumnWidth = 12.7109375
$ws.Columns.Item(5).ColumnWidth = 12.5703125
$ws.Columns.Item(7).ColumnWidth = 32.85546875
$ws.Columns.Item(8).ColumnWidth = 21.5703125
$ws.Columns.Item(9).ColumnWidth = 83.42578125
